# Update Name of Algo
# Applies updated RandomForest imputation values to specific cells in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -8.7043
$ws.Range("E3").Value = 15.80940000000001
$ws.Range("A12").Value = -21.5522
$ws.Range("D14").Value = -7.4073
$ws.Range("E20").Value = 15.8962
$ws.Range("E25").Value = 17.20910000000001
$ws.Range("D26").Value = -8.620799999999996
$ws.Range("A27").Value = -21.59229999999999
$ws.Range("E30").Value = 15.5257
$ws.Range("D31").Value = -8.734199999999998
$ws.Range("A32").Value = -21.3351
$ws.Range("D35").Value = -8.938599999999994
$ws.Range("A36").Value = -20.07660000000001
$ws.Range("D37").Value = -8.606499999999995
$ws.Range("A38").Value = -19.59070000000001
$ws.Range("E44").Value = 16.5428
$ws.Range("D45").Value = -7.568700000000002
$ws.Range("A46").Value = -21.689
$ws.Range("E47").Value = 16.00369999999999
$ws.Range("D52").Value = -7.868999999999999
$ws.Range("A54").Value = -21.5914
$ws.Range("A55").Value = -22.49350000000002
$ws.Range("A56").Value = -22.22120000000001
$ws.Range("D57").Value = -8.628199999999994
$ws.Range("E58").Value = 16.64040000000001
$ws.Range("A67").Value = -21.47699999999997
$ws.Range("A69").Value = -21.60209999999998
$ws.Range("A72").Value = -21.75710000000001
$ws.Range("E78").Value = 16.62460000000003
$ws.Range("D81").Value = -7.147499999999995
$ws.Range("A83").Value = -21.9237
$ws.Range("D83").Value = -8.315900000000001
$ws.Range("E84").Value = 16.3936
$ws.Range("A86").Value = -21.97670000000002
$ws.Range("E89").Value = 17.34630000000002
$ws.Range("A91").Value = -21.4681
$ws.Range("E91").Value = 17.91470000000002
$ws.Range("E92").Value = 17.98540000000002
$ws.Range("A93").Value = -21.23369999999999
$ws.Range("E96").Value = 15.21909999999999
$ws.Range("A99").Value = -19.7128
$ws.Range("D100").Value = -8.367999999999999
$ws.Range("D102").Value = -7.7331
$ws.Range("E102").Value = 16.57860000000001
